$wb = $excel.ActiveWorkbook

# The "想去人数" (number of people wanting to go) counts were updated for two
# events, on both the "展览" sheet and the mirrored "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 261
    $ws.Range("F4").Value = 162
}
